# ---------------------------------------------------------------------------
# Add prepared results for the "eff 10p down" sensitivity run to the
# Comparison_sensitivities sheet (Table3), wire it into the LCOE sensitivity
# chart as a new "Efficiency" series, and make that sheet the active tab.
# ---------------------------------------------------------------------------

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Comparison_operating_points")
$ws3 = $wb.Worksheets.Item("Comparison_sensitivities")

# ---------------------------------------------------------------------------
# 1. New result rows (63 & 64) for run "sens_10op_eff_10pdown" (+ its PV
#    variant) at the bottom of the Table3 data range.
# ---------------------------------------------------------------------------

$ws3.Range("A63").Value = "sens_10op_eff_10pdown"
$ws3.Range("B63").Formula = '=IF(ISNUMBER(SEARCH("PV", A63)),"PV revenue","no PV revenue")'
$ws3.Range("C63").Value = 262.62631520402192
$ws3.Range("D63").Value = 1451.739909044454
$ws3.Range("E63").Value = 298331833.89240187
$ws3.Range("F63").Value = 16069931.11681715
$ws3.Range("G63").Value = 176888.92799996189
$ws3.Range("H63").Value = 32000.007075370009
$ws3.Range("I63").Value = 9.8181474074492936

$ws3.Range("A64").Value = "sens_10op_eff_10pdown_PV"
$ws3.Range("B64").Formula = '=IF(ISNUMBER(SEARCH("PV", A64)),"PV revenue","no PV revenue")'
$ws3.Range("C64").Value = 211.80775199938211
$ws3.Range("D64").Value = 1170.8261846632511
$ws3.Range("E64").Value = 298331833.89240187
$ws3.Range("F64").Value = 7080689.949050107
$ws3.Range("G64").Value = 176888.92799996189
$ws3.Range("H64").Value = 32000.007075370009
$ws3.Range("I64").Value = 9.8181474074492936

# Grow Table3 (run_name .. pcf_value) so the autofilter/table range covers
# the two new rows.
$lo3 = $ws3.ListObjects.Item("Table3")
$lo3.Resize($ws3.Range("A1:I64"))

# ---------------------------------------------------------------------------
# 2. New "Efficiency" sensitivity summary column (K..S helper block) so the
#    chart has a source range to plot.
# ---------------------------------------------------------------------------

$ws3.Range("S1").Value = "Efficiency"
$ws3.Range("S2").Formula = "=D63"
$ws3.Range("S3").Value = 1400
$ws3.Range("S4").Formula = "=`$D`$2"

# Cross-check columns (were T4:U4 / T5:U5, now shifted one column right to
# U4:V4 / U5:V5 to make room for the new "Efficiency" column S).
$ws3.Range("U4").Formula = "=R2/P4"
$ws3.Range("V4").Formula = "=1-U4"
$ws3.Range("U5").Formula = "=R6/P4"
$ws3.Range("V5").Formula = "=U5-1"
$ws3.Range("T4").Clear()
$ws3.Range("T5").Clear()

# ---------------------------------------------------------------------------
# 3. Chart: fix the "InvCost" series (was missing its category axis range)
#    and add the new "Efficiency" series plotted the same way.
# ---------------------------------------------------------------------------

$chart = $ws3.ChartObjects(1).Chart
$sc = $chart.SeriesCollection()

$invCost = $sc.Item(7)
$invCost.Formula = "=SERIES(Comparison_sensitivities!`$R`$1,Comparison_sensitivities!`$K`$2:`$K`$6,Comparison_sensitivities!`$R`$2:`$R`$6,7)"

$effSeries = $sc.NewSeries()
$effSeries.Formula = "=SERIES(Comparison_sensitivities!`$S`$1,Comparison_sensitivities!`$K`$2:`$K`$6,Comparison_sensitivities!`$S`$2:`$S`$6,8)"
$effSeries.Name = "Efficiency"

# ---------------------------------------------------------------------------
# 4. Cosmetic cleanup to match the re-saved workbook: the helper "PV/no PV"
#    column on the operating-points sheet loses its explicit (no-op) number
#    format, and Comparison_sensitivities becomes the active tab.
# ---------------------------------------------------------------------------

$ws1.Range("B22:B25").Style = "Normal"

$ws3.Activate()
$ws3.Range("S4").Select()
